$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix tiny floating point precision update on A10 (recomputed timestamp)
$ws.Range("A10").Value = 45862.79190449074

# Append new row 11 with latest sensor reading
$ws.Range("A11").Value = 45862.87522489934
$ws.Range("A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value = 2025
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = 19.33
$ws.Range("E11").Value = 74.44
$ws.Range("F11").Value = 87.89
$ws.Range("G11").Value = 13.54
$ws.Range("H11").Value = "ESE"
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "21:00:19"
